$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = 'Alyne Corona_20251128_162805'
$ws.Range("B17").Value = "'"
$ws.Range("B17").ClearFormats()
$ws.Range("C17").Value = 'Alyne Corona'
$ws.Range("D17").Value = 21
$ws.Range("E17").Value = 'Female'
$ws.Range("F17").Value = '2025-11-28 16:28:05'
$ws.Range("G17").Value = '{
  "portion": 0.4,
  "diet": 0.2857142857142857,
  "salt": 0.2,
  "fat": 0.2,
  "natural": 1.0,
  "convenience": 0.2,
  "price": 0.2
}'
$ws.Rows.Item(17).AutoFit()
$ws.Range("H17").Value = 'Nongshim Neoguri Spicy Seafood'
$ws.Range("I17").NumberFormat = "@"
$ws.Range("I17").Value = '0.723'
$ws.Range("I17").ClearFormats()
$ws.Range("J17").Value = 'Sabor a marisco, umami, picante equilibrado, buena textura, algo salado'
$ws.Range("K17").Value = 'Nissin Chow Mein Teriyaki Beef'
$ws.Range("L17").NumberFormat = "@"
$ws.Range("L17").Value = '0.370'
$ws.Range("L17").ClearFormats()
$ws.Range("M17").Value = 'Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa'
$ws.Range("N17").Value = 'Nongshim Shin Ramyun'
$ws.Range("O17").NumberFormat = "@"
$ws.Range("O17").Value = '0.337'
$ws.Range("O17").ClearFormats()
$ws.Range("P17").Value = 'Sabor intenso, picante, umami, fideos gruesos, muy alto en sodio'
$ws.Range("Q17").Value = 'Amy’s Macaroni & Cheese (frozen)'
$ws.Range("R17").NumberFormat = "@"
$ws.Range("R17").Value = '0.758'
$ws.Range("R17").ClearFormats()
$ws.Range("S17").Value = 'Queso real, textura casera, sin conservadores, alto en grasa, algo caro'
$ws.Range("T17").Value = 'Kraft Macaroni & Cheese Dinner'
$ws.Range("U17").NumberFormat = "@"
$ws.Range("U17").Value = '0.715'
$ws.Range("U17").ClearFormats()
$ws.Range("V17").Value = 'Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato'
$ws.Range("W17").Value = 'Annie’s Shells & White Cheddar'
$ws.Range("X17").NumberFormat = "@"
$ws.Range("X17").Value = '0.686'
$ws.Range("X17").ClearFormats()
$ws.Range("Y17").Value = 'Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños'
$ws.Range("Z17").Value = 'Wild Planet Wild Tuna Pasta Salad'
$ws.Range("AA17").NumberFormat = "@"
$ws.Range("AA17").Value = '0.795'
$ws.Range("AA17").ClearFormats()
$ws.Range("AB17").Value = 'Sabor fresco, buena proteína, saludable, porción algo pequeña'
$ws.Range("AC17").Value = 'Kitchens of India Variety Pack'
$ws.Range("AD17").NumberFormat = "@"
$ws.Range("AD17").Value = '0.653'
$ws.Range("AD17").ClearFormats()
$ws.Range("AE17").Value = 'Sabor auténtico, variedad, vegetariano, necesita arroz o pan, buena calidad'
$ws.Range("AF17").Value = 'StarKist Chicken Creations (Chicken Salad)'
$ws.Range("AG17").NumberFormat = "@"
$ws.Range("AG17").Value = '0.359'
$ws.Range("AG17").ClearFormats()
$ws.Range("AH17").Value = 'Portátil, saludable, fácil, buena textura, sabor suave'
